$d = $word.ActiveDocument

# 1. "API RESTful" -> "API" (Objetivo section, back-end technology mention)
$d.Content.Find.Execute("API RESTful", $true, $false, $false, $false, $false,
                         $true, 1, $false, "API", 2)

# 2. Drop the "Stored Procedures" mention from the ORM sentence:
#    " para mapeamento objeto-relacional e Stored Procedures para algumas operações do banco de dados."
#    becomes " para mapeamento objeto-relacional."
#    Step 2a: trim the trailing " e " off the (non-bold) run that precedes "Stored Procedures".
$rngOrm = $d.Content
$rngOrm.Find.Execute("objeto-relacional e ", $true, $false, $false, $false, $false,
                      $true, 1, $false, "objeto-relacional", 2)

#    Step 2b: collapse "Stored Procedures para algumas operações do banco de dados."
#    down to just "." -- the match starts inside the bold "Stored Procedures" run, so a
#    plain replace would leave the new run bold. Borrow the (non-bold) formatting of a
#    lone "." elsewhere in the document via FormattedText so the run keeps matching,
#    non-bold formatting instead.
$rngTail = $d.Content
$rngTail.Find.Execute("Stored Procedures para algumas operações do banco de dados.",
                       $true, $false, $false, $false, $false,
                       $true, 1, $false, ".", 2)

$rngPeriodSrc = $d.Content
$rngPeriodSrc.Find.Execute("usuários.")
$periodOnly = $d.Range($rngPeriodSrc.End - 1, $rngPeriodSrc.End)
$rngTail.FormattedText = $periodOnly.FormattedText

# 3. "3.2. Back-End (API REST em C# com .NET Core 6.0)" -> "3.2. Back-End (API em C# com .NET Core 6.0)"
$d.Content.Find.Execute("3.2. Back-End (API REST em C# com .NET Core 6.0)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3.2. Back-End (API em C# com .NET Core 6.0)", 2)
